# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Mon Aug  7 20:51:01 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values look like plain numbers (e.g. "241.24", "1.000").
# Force just those specific cells to Text first so Excel stores the digits
# verbatim instead of silently parsing them into a number.
$textForcedRows = @(4,5,6,7,8,9,10,11,13,14,15,16,17,20,21,22,23,24,25,26,27,28,30,31,32,33,34,35,36,37,38,41,42,43,44,46,47,48,49,50,51)
foreach ($r in $textForcedRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "29.116.12"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3
$ws.Range("D3").Value = "1.822.76"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").Value = "241.24"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").Value = "0.6153"
$ws.Range("E6").Value = "  -1.99%  "

# Row 7
$ws.Range("D7").Value = "0.9988"
$ws.Range("E7").Value = "  -0.29%  "

# Row 8
$ws.Range("D8").Value = "0.07331"
$ws.Range("E8").Value = "  -2.01%  "

# Row 9
$ws.Range("D9").Value = "0.2890"
$ws.Range("E9").Value = "  -1.08%  "

# Row 10
$ws.Range("D10").Value = "22.94"
$ws.Range("E10").Value = "  -0.92%  "

# Row 11
$ws.Range("D11").Value = "0.07656"
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("D12").Value = "1.826.94"
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").Value = "4.948"
$ws.Range("E13").Value = "  -1.08%  "

# Row 14
$ws.Range("D14").Value = "0.6600"
$ws.Range("E14").Value = "  -1.10%  "

# Row 15
$ws.Range("D15").Value = "81.79"
$ws.Range("E15").Value = "  -1.15%  "

# Row 16
$ws.Range("D16").Value = "0.000008936"
$ws.Range("E16").Value = "  -5.00%  "

# Row 17
$ws.Range("D17").Value = "5.828"
$ws.Range("E17").Value = "  -2.43%  "

# Row 18
$ws.Range("D18").Value = "29.072.83"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
$ws.Range("D19").Value = "2.069.24"
$ws.Range("E19").Value = "  +0.00%  "

# Row 20
$ws.Range("D20").Value = "236.86"
$ws.Range("E20").Value = "  +6.43%  "

# Row 21
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.34%  "

# Row 23
$ws.Range("D23").Value = "7.119"
$ws.Range("E23").Value = "  +0.26%  "

# Row 24
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("D25").Value = "157.62"
$ws.Range("E25").Value = "  -1.35%  "

# Row 26
$ws.Range("D26").Value = "0.1412"
$ws.Range("E26").Value = "  +1.54%  "

# Row 27
$ws.Range("D27").Value = "8.423"
$ws.Range("E27").Value = "  -0.68%  "

# Row 28
$ws.Range("D28").Value = "17.61"
$ws.Range("E28").Value = "  -1.36%  "

# Row 29
$ws.Range("E29").Value = "  -0.98%  "

# Row 30
$ws.Range("D30").Value = "0.05541"
$ws.Range("E30").Value = "  -3.49%  "

# Row 31
$ws.Range("D31").Value = "4.082"
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("D32").Value = "4.090"
$ws.Range("E32").Value = "  -1.32%  "

# Row 33
$ws.Range("D33").Value = "1.206"
$ws.Range("E33").Value = "  -0.20%  "

# Row 34
$ws.Range("D34").Value = "1.823"
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
$ws.Range("D35").Value = "0.7334"
$ws.Range("E35").Value = "  -0.71%  "

# Row 36
$ws.Range("D36").Value = "1.131"
$ws.Range("E36").Value = "  -0.52%  "

# Row 37
$ws.Range("D37").Value = "2.608"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38
$ws.Range("D38").Value = "2.833"
$ws.Range("E38").Value = "  +2.40%  "

# Row 39
$ws.Range("D39").Value = "1.206.58"
$ws.Range("E39").Value = "  -0.67%  "

# Row 40
$ws.Range("E40").Value = "  -1.25%  "

# Row 41
$ws.Range("D41").Value = "6.331"
$ws.Range("E41").Value = "  -2.46%  "

# Row 42
$ws.Range("D42").Value = "0.8979"
$ws.Range("E42").Value = "  +0.90%  "

# Row 43
$ws.Range("D43").Value = "0.9987"
$ws.Range("E43").Value = "  -0.29%  "

# Row 44
$ws.Range("D44").Value = "101.11"
$ws.Range("E44").Value = "  -0.96%  "

# Row 45
$ws.Range("D45").Value = "1.975.30"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  +0.35%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "64.50"
$ws.Range("E47").Value = "  -1.51%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.5079"
$ws.Range("E48").Value = "  -0.21%  "

# Row 49
$ws.Range("D49").Value = "0.3997"
$ws.Range("E49").Value = "  -1.49%  "

# Row 50
$ws.Range("D50").Value = "8.991"
$ws.Range("E50").Value = "  +0.36%  "

# Row 51
$ws.Range("D51").Value = "0.05748"
$ws.Range("E51").Value = "  -1.27%  "
